$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; this shifts all existing data rows
# down by one (rows 1-15 -> 2-16) and the summary row (25 -> 26), and
# Excel auto-updates the GEOMEAN() ranges in the summary formulas.
$ws.Rows("1").Insert()

# Fill in the new header row. Order matches the order the strings were
# first introduced in the authoritative edit (f/s accel, s/f accel,
# P_their, P_our, T_their, T_our) so the shared-strings table comes out
# in the same sequence.
$ws.Range("C1").Value = "f/s acceleration"
$ws.Range("H1").Value = "s/f acceleration"
$ws.Range("F1").Value = "P_their"
$ws.Range("G1").Value = "P_our"
$ws.Range("A1").Value = "T_their"
$ws.Range("B1").Value = "T_our"

# Re-point the two formula columns cell-by-cell (rather than as one
# multi-cell range write) so each row keeps a formula that resolves to
# the correct relative references after the shift.
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Formula = "=A$r/B$r"
    $ws.Cells.Item($r, 8).Formula = "=G$r/F$r"
}

# Match the author's final selection.
$ws.Range("G1").Select() | Out-Null
